$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A7").Value = "Fanta"

$ws.Range("B7").Value = "'20"
$ws.Range("B7").ClearFormats()

$ws.Range("C7").Value = "'50"
$ws.Range("C7").ClearFormats()
